$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the student-name placeholders in column A (rows 2-6) with numeric IDs,
# dropping the "NaN" (non-numeric) inline string values in favor of numbers.
$ws.Range("A2").Value = 11
$ws.Range("A3").Value = 12
$ws.Range("A4").Value = 13
$ws.Range("A5").Value = 14
$ws.Range("A6").Value = 15
